# STAG Notes.docx robustness edit:
#  - remove the stray empty bullet under "Health"
#  - bold the lead keyword of each "Action" bullet (Trigger / Subject /
#    Consumed / Produced) and split it into its own run
#  - add "e.g." sub-bullets under each of those bullets
#  - add two new bullets describing the action-verification workflow

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Delete the empty ilvl=2 bullet paragraph that sits right after the
#    "Health" bullet (and right before "Actions - ...").
# ---------------------------------------------------------------------
$healthRng = $d.Content
$healthFound = $healthRng.Find.Execute("Health", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($healthFound) {
    $healthPara = $healthRng.Paragraphs(1)
    $emptyPara = $healthPara.Next()
    if ($emptyPara.Range.Text.Trim().Length -eq 0) {
        $emptyPara.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# Helper: find a paragraph by its (unique) leading text, make the first
# $boldLen characters bold (splitting the run), and return the Paragraph.
# ---------------------------------------------------------------------
function Bold-Lead($searchText, $boldLen) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    $lead = $d.Range($rng.Start, $rng.Start + $boldLen)
    $lead.Font.Bold = 1
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------
# Helper: insert a brand-new bullet paragraph right after $afterPara,
# at outline level $level (1 = top, 2 = second, 3 = third / ilvl 2),
# containing $text.
# ---------------------------------------------------------------------
function Insert-BulletAfter($afterPara, $text, $level) {
    $afterPara.Range.InsertParagraphAfter()
    $newPara = $afterPara.Next()
    $newPara.Range.Text = $text
    $newPara.Range.ListFormat.ListLevelNumber = $level
    return $newPara
}

# ---------------------------------------------------------------------
# 2. "Trigger words - Initiate the action."
# ---------------------------------------------------------------------
$triggerPara = Bold-Lead "Trigger words" 7
if ($triggerPara -ne $null) {
    Insert-BulletAfter $triggerPara "e.g., open, unlock." 3 | Out-Null
}

# ---------------------------------------------------------------------
# 3. "Subject entities - Entities that must be present in inv and
#     location to do action."
# ---------------------------------------------------------------------
$subjectPara = Bold-Lead "Subject entities" 7
if ($subjectPara -ne $null) {
    Insert-BulletAfter $subjectPara "E.g., door, key." 3 | Out-Null
}

# ---------------------------------------------------------------------
# 4. "Consumed entities - Entities that are removed/eaten up by action."
# ---------------------------------------------------------------------
$consumedPara = Bold-Lead "Consumed entities" 9
if ($consumedPara -ne $null) {
    Insert-BulletAfter $consumedPara "E.g., key" 3 | Out-Null
}

# ---------------------------------------------------------------------
# 5. "Produced entities - Entities created as a result of action."
# ---------------------------------------------------------------------
$producedPara = Bold-Lead "Produced entities" 9
if ($producedPara -ne $null) {
    $cellarPara = Insert-BulletAfter $producedPara "E.g., cellar" 3

    # ---------------------------------------------------------------
    # 6. Two brand-new bullets describing the verification workflow.
    # ---------------------------------------------------------------
    $enginePara = Insert-BulletAfter $cellarPara "Game engine should verify conditions for actions to be performed," 2
    Insert-BulletAfter $enginePara "Then undertake relevant additions/removals" 3 | Out-Null
}

Write-Output "done"
